$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Table on slide 5: switch the table style (tableStyleId) to the new
#    built-in style GUID.
# ---------------------------------------------------------------------------
$tableSlide = $p.Slides.Item(5)
for ($i = 1; $i -le $tableSlide.Shapes.Count; $i++) {
    $shp = $tableSlide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{E7E0E156-9413-4380-9463-799A5DFCDE02}")
    }
}

# ---------------------------------------------------------------------------
# 2) Re-colour the presentation's theme (the "Integral" / Red Violet scheme)
#    to the standard Office colour scheme (dk1/lt1/dk2/lt2/accent1-6/
#    hlink/folHlink), matching the swapped theme content.
# ---------------------------------------------------------------------------
$master = $p.SlideMaster
$colorScheme = $master.Theme.ThemeColorScheme

# Index order exposed by ThemeColorScheme: dk1, lt1, dk2, lt2,
# accent1, accent2, accent3, accent4, accent5, accent6, hlink, folHlink.
$officeColors = @(0, 16777215, 6968388, 15132391, 13998939, 3243501, 10855845, 49407, 12874308, 4697456, 12673797, 7491477)

for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $colorScheme.Item($i).RGB = $officeColors[$i - 1]
}
